$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.13541
$ws.Range("H2").Value = 0.40623
$ws.Range("I2").Value = 0.03919062573893041
$ws.Range("J2").Value = 0.0391906257389304
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.745085
$ws.Range("N2").Value = 2.235255
$ws.Range("O2").Value = 0.07418358086139912
$ws.Range("P2").Value = 0.07418358086139913
$ws.Range("Q2").Value = 0.10089195985
$ws.Range("R2").Value = 0.9080276386499999
$ws.Range("S2").Value = 0.002907300953512774
$ws.Range("T2").Value = 0.002907300953512774

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.13541
$ws.Range("H3").Value = 0.40623
$ws.Range("I3").Value = 0.03919062573893041
$ws.Range("J3").Value = 0.0391906257389304
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.332302333333334
$ws.Range("N3").Value = 15.996907
$ws.Range("O3").Value = 0.5309049052420336
$ws.Range("P3").Value = 0.5309049052420336
$ws.Range("Q3").Value = 0.7220470589566668
$ws.Range("R3").Value = 6.498423530609999
$ws.Range("S3").Value = 0.02080649544430285
$ws.Range("T3").Value = 0.02080649544430285

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra2"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.13541
$ws.Range("H4").Value = 0.40623
$ws.Range("I4").Value = 0.03919062573893041
$ws.Range("J4").Value = 0.0391906257389304
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07052433333333334
$ws.Range("N4").Value = 0.211573
$ws.Range("O4").Value = 0.00702167884809062
$ws.Range("P4").Value = 0.007021678848090619
$ws.Range("Q4").Value = 0.009549699976666668
$ws.Range("R4").Value = 0.08594729979
$ws.Range("S4").Value = 0.0002751839877944835
$ws.Range("T4").Value = 0.0002751839877944834

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra2"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.13541
$ws.Range("H5").Value = 0.40623
$ws.Range("I5").Value = 0.03919062573893041
$ws.Range("J5").Value = 0.0391906257389304
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.895887666666667
$ws.Range("N5").Value = 11.687663
$ws.Range("O5").Value = 0.3878898350484767
$ws.Range("P5").Value = 0.3878898350484767
$ws.Range("Q5").Value = 0.5275421489433334
$ws.Range("R5").Value = 4.74787934049
$ws.Range("S5").Value = 0.0152016453533203
$ws.Range("T5").Value = 0.0152016453533203

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Gfra2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.319753
$ws.Range("H6").Value = 9.959258999999999
$ws.Range("I6").Value = 0.9608093742610696
$ws.Range("J6").Value = 0.9608093742610695
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.745085
$ws.Range("N6").Value = 2.235255
$ws.Range("O6").Value = 0.07418358086139912
$ws.Range("P6").Value = 0.07418358086139913
$ws.Range("Q6").Value = 2.473498164005
$ws.Range("R6").Value = 22.261483476045
$ws.Range("S6").Value = 0.07127627990788635
$ws.Range("T6").Value = 0.07127627990788635

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Gfra2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.319753
$ws.Range("H7").Value = 9.959258999999999
$ws.Range("I7").Value = 0.9608093742610696
$ws.Range("J7").Value = 0.9608093742610695
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.332302333333334
$ws.Range("N7").Value = 15.996907
$ws.Range("O7").Value = 0.5309049052420336
$ws.Range("P7").Value = 0.5309049052420336
$ws.Range("Q7").Value = 17.70192666799034
$ws.Range("R7").Value = 159.317340011913
$ws.Range("S7").Value = 0.5100984097977308
$ws.Range("T7").Value = 0.5100984097977307

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Gdnf"
$ws.Range("C8").Value = "Gfra2"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.319753
$ws.Range("H8").Value = 9.959258999999999
$ws.Range("I8").Value = 0.9608093742610696
$ws.Range("J8").Value = 0.9608093742610695
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07052433333333334
$ws.Range("N8").Value = 0.211573
$ws.Range("O8").Value = 0.00702167884809062
$ws.Range("P8").Value = 0.007021678848090619
$ws.Range("Q8").Value = 0.2341233671563334
$ws.Range("R8").Value = 2.107110304407
$ws.Range("S8").Value = 0.006746494860296137
$ws.Range("T8").Value = 0.006746494860296135

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Gdnf"
$ws.Range("C9").Value = "Gfra2"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.319753
$ws.Range("H9").Value = 9.959258999999999
$ws.Range("I9").Value = 0.9608093742610696
$ws.Range("J9").Value = 0.9608093742610695
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.895887666666667
$ws.Range("N9").Value = 11.687663
$ws.Range("O9").Value = 0.3878898350484767
$ws.Range("P9").Value = 0.3878898350484767
$ws.Range("Q9").Value = 12.93338476907967
$ws.Range("R9").Value = 116.400462921717
$ws.Range("S9").Value = 0.3726881896951564
$ws.Range("T9").Value = 0.3726881896951563

# Remove now-obsolete rows 10 and 11 (MuSCs -> MuSCs/Resolving-Mac duplicate rows)
$ws.Rows("10:11").Delete()
